$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per the latest crypto data pull.
# D values are entered with a leading apostrophe so Excel keeps them as text
# (many look numeric, e.g. "181.64", and would otherwise be auto-converted to
# a Number cell). The Style reset afterwards clears the resulting quote-prefix
# formatting so the cell keeps the workbook default style, matching the source data.
$ws.Range("D2").Value = "'65.533.57"
$ws.Range("E2").Value = "  -6.23%  "
$ws.Range("D3").Value = "'3.286.38"
$ws.Range("E3").Value = "  -6.51%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'553.12"
$ws.Range("E5").Value = "  -4.62%  "
$ws.Range("D6").Value = "'181.64"
$ws.Range("E6").Value = "  -5.92%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.586"
$ws.Range("E8").Value = "  -4.48%  "
$ws.Range("D9").Value = "'3.281.88"
$ws.Range("E9").Value = "  -6.34%  "
$ws.Range("D10").Value = "'0.183"
$ws.Range("E10").Value = "  -10.78%  "
$ws.Range("D11").Value = "'0.581"
$ws.Range("E11").Value = "  -6.69%  "
$ws.Range("D12").Value = "'47.05"
$ws.Range("E12").Value = "  -8.69%  "
$ws.Range("D13").Value = "'0.0000265"
$ws.Range("E13").Value = "  -7.70%  "
$ws.Range("D14").Value = "'642.53"
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("D15").Value = "'8.62"
$ws.Range("E15").Value = "  -6.23%  "
$ws.Range("D16").Value = "'3.817.72"
$ws.Range("E16").Value = "  -6.27%  "
$ws.Range("D17").Value = "'18.03"
$ws.Range("E17").Value = "  -2.10%  "
$ws.Range("D18").Value = "'65.607.38"
$ws.Range("E18").Value = "  -6.08%  "
$ws.Range("E19").Value = "  -3.27%  "
$ws.Range("D20").Value = "'3.301.90"
$ws.Range("E20").Value = "  -6.27%  "
$ws.Range("D21").Value = "'11.34"
$ws.Range("E21").Value = "  -8.93%  "
$ws.Range("D22").Value = "'0.901"
$ws.Range("E22").Value = "  -5.49%  "
$ws.Range("D23").Value = "'17.98"
$ws.Range("E23").Value = "  -0.84%  "
$ws.Range("D24").Value = "'107.12"
$ws.Range("E24").Value = "  +7.46%  "
$ws.Range("D25").Value = "'4.87"
$ws.Range("E25").Value = "  -9.14%  "
$ws.Range("D26").Value = "'3.95"
$ws.Range("E26").Value = "  -8.27%  "
$ws.Range("D27").Value = "'2.67"
$ws.Range("E27").Value = "  -8.08%  "
$ws.Range("D28").Value = "'9.51"
$ws.Range("E28").Value = "  -6.29%  "
$ws.Range("D29").Value = "'8.64"
$ws.Range("E29").Value = "  -8.39%  "
$ws.Range("D30").Value = "'30.08"
$ws.Range("E30").Value = "  -8.37%  "
$ws.Range("D31").Value = "'3.87"
$ws.Range("E31").Value = "  -9.75%  "
$ws.Range("D32").Value = "'6.23"
$ws.Range("E32").Value = "  -8.12%  "
$ws.Range("D33").Value = "'10.99"
$ws.Range("E33").Value = "  -5.87%  "
$ws.Range("E34").Value = "  -5.54%  "
$ws.Range("D35").Value = "'3.771.34"
$ws.Range("E35").Value = "  -0.24%  "
$ws.Range("D36").Value = "'57.32"
$ws.Range("E36").Value = "  -6.76%  "
$ws.Range("D38").Value = "'520.17"
$ws.Range("E38").Value = "  -9.50%  "
$ws.Range("D39").Value = "'0.0₃0726"
$ws.Range("E39").Value = "  -8.84%  "
$ws.Range("D40").Value = "'3.32"
$ws.Range("E40").Value = "  -8.19%  "
$ws.Range("D41").Value = "'0.128"
$ws.Range("E41").Value = "  -3.69%  "
$ws.Range("D42").Value = "'2.69"
$ws.Range("E42").Value = "  -7.25%  "
$ws.Range("D43").Value = "'32.82"
$ws.Range("E43").Value = "  -4.46%  "
$ws.Range("D44").Value = "'3.29"
$ws.Range("E44").Value = "  -12.57%  "
$ws.Range("D45").Value = "'0.335"
$ws.Range("E45").Value = "  -10.82%  "
$ws.Range("E46").Value = "  -2.52%  "
$ws.Range("D47").Value = "'0.0412"
$ws.Range("E47").Value = "  -7.54%  "
$ws.Range("E48").Value = "  -5.15%  "
$ws.Range("D49").Value = "'2.59"
$ws.Range("E49").Value = "  -10.39%  "
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("E51").Value = "  +1.26%  "

# Clear the quote-prefix style artifact introduced above so D-column cells
# retain the workbook default (unstyled) appearance.
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
